# unify the conception of DataNode, DataTable, Entity.
#
# Renames the "Property*" sheets to "DataNode_*" and the "Record_*" sheets
# to "DataTable_*", drops the obsolete "Record_Building" sheet (and its
# building-specific rows/comments go with it), and makes "DataTable_Hero"
# (formerly "Record_Hero") the selected/active sheet.

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

$wb.Worksheets.Item("Property1").Name = "DataNode_1"
$wb.Worksheets.Item("Property2").Name = "DataNode_2"
$wb.Worksheets.Item("Record_Hero").Name = "DataTable_Hero"
$wb.Worksheets.Item("Record_Bag").Name = "DataTable_Bag"
$wb.Worksheets.Item("Record_CommPropertyValue").Name = "DataTable_CommPropertyValue"
$wb.Worksheets.Item("Record_Task").Name = "DataTable_Task"

# The building-record table is no longer part of the unified schema.
$wb.Worksheets.Item("Record_Building").Delete()

# DataTable_Hero becomes the active/selected sheet.
$wb.Worksheets.Item("DataTable_Hero").Activate()
